$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 3980
$ws1.Range("F4").Value = 2336
$ws1.Range("F9").Value = 195
$ws1.Range("F11").Value = 53
$ws1.Range("F12").Value = 122
$ws1.Range("F13").Value = 1482
$ws1.Range("F15").Value = 2752
$ws1.Range("F16").Value = 190

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 3980
$ws4.Range("F4").Value = 2336
$ws4.Range("F10").Value = 195
$ws4.Range("F12").Value = 53
$ws4.Range("F13").Value = 122
$ws4.Range("F16").Value = 1482
$ws4.Range("F18").Value = 2752
$ws4.Range("F19").Value = 190
